$d = $word.ActiveDocument

# Locate the anchor paragraph: the "Requisitos funcionais" bullet (ilvl=1, numId=4)
# that immediately precedes the "Requisitos de qualidade" bullet. (There is an
# earlier, unrelated "Requisitos funcionais " heading higher up in the outline,
# so matching on adjacency to "Requisitos de qualidade" disambiguates the two.)
$anchor = $null
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $next = $d.Paragraphs.Item($i + 1)
    if (($p.Range.Text.Trim() -eq "Requisitos funcionais") -and ($next.Range.Text.Trim() -eq "Requisitos de qualidade")) {
        $anchor = $p
    }
}
if ($anchor -eq $null) {
    # Fallback: last paragraph whose trimmed text is exactly "Requisitos funcionais".
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq "Requisitos funcionais") {
            $anchor = $p
        }
    }
}

$normalStyle = $d.Styles.Item("Normal")

# (ListLevelNumber, text) pairs to insert, in document order, right after the anchor.
$newItems = @(
    @(3, "Assim que iniciado, o sistema deverá mostrar ao usuário um menu e verificar que operação ele deseja realizar:"),
    @(4, " verificar se um número é par ou ímpar;"),
    @(4, "verificar se um numero é positivo ou negativo;"),
    @(4, "verificar se ele deseja sair do programa."),
    @(3, "O sistema deverá verificar qual a opção que o usuário deseja realizar."),
    @(3, "O sistema deverá verificar se um número é ímpar ou par quando um número inteiro positivo é fornecido"),
    @(3, "O sistema deverá informar ao usuário se o número é negativo ou positivo quando um número inteiro é fornecido pelo usuário"),
    @(3, "O sistema deverá exibir uma mensagem para notificar o usuário que o sistema foi finalizado")
)

$prev = $anchor
foreach ($item in $newItems) {
    $level = $item[0]
    $text = $item[1]

    $prev.Range.InsertParagraphAfter()
    $newPara = $prev.Next()

    $newPara.Range.Text = $text
    $newPara.Style = $normalStyle
    $newPara.Range.ListFormat.ListLevelNumber = $level

    $prev = $newPara
}

Write-Output "Inserted $($newItems.Count) paragraphs after 'Requisitos funcionais'."
